$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 461.5
$ws.Range("I8").Value = 461.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1384.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1245.5

$ws.Range("H38").Value = 228.18182
$ws.Range("I38").Value = 158
$ws.Range("J38").Value = 544
$ws.Range("K38").Value = 474
$ws.Range("L38").Value = 1632
$ws.Range("M38").Value = -102
$ws.Range("N38").Value = -2376

$ws.Range("H62").Value = 3629.6667
$ws.Range("I62").Value = 4444.5
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 4444.5
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -3820.5
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 3629.6667
$ws.Range("I65").Value = 4444.5
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 22222.5
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -19102.5
$ws.Range("N65").Value = -16240

$ws.Range("H100").Value = 2999.5
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 2999
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 2999
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -4081

$ws.Range("H116").Value = 3047.7368
$ws.Range("I116").Value = 2554.0908
$ws.Range("J116").Value = 3726.5
$ws.Range("K116").Value = 2554.0908
$ws.Range("L116").Value = 3726.5
$ws.Range("M116").Value = 887.9092000000001
$ws.Range("N116").Value = -10610.5

$ws.Range("H126").Value = 30000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -39880

$ws.Range("H138").Value = 2127.97
$ws.Range("I138").Value = 1633.2222
$ws.Range("J138").Value = 2176.9011
$ws.Range("K138").Value = 4899.6666
$ws.Range("L138").Value = 6530.7033
$ws.Range("M138").Value = 240.3334000000004
$ws.Range("N138").Value = -16810.7033

$ws.Range("H141").Value = 7564.6665
$ws.Range("I141").Value = 9027.308000000001
$ws.Range("J141").Value = 3761.8
$ws.Range("K141").Value = 27081.924
$ws.Range("L141").Value = 11285.4
$ws.Range("M141").Value = -21901.924
$ws.Range("N141").Value = -21645.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5006.4
$ws.Range("I32").Value = 5006.4
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5006.4
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4719.4
$ws.Range("N32").Value = ""

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""

$ws.Range("H61").Value = 950.2632
$ws.Range("I61").Value = 776.6875
$ws.Range("J61").Value = 1876
$ws.Range("K61").Value = 776.6875
$ws.Range("L61").Value = 1876
$ws.Range("M61").Value = -564.6875
$ws.Range("N61").Value = -2300

$ws.Range("H63").Value = 2950
$ws.Range("I63").Value = 2900
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2900
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -2214
$ws.Range("N63").Value = -4372

$ws.Range("H66").Value = 2950
$ws.Range("I66").Value = 2900
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 14500
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -11068
$ws.Range("N66").Value = -21864

$ws.Range("H80").Value = 37000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 37000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 37000
$ws.Range("N80").Value = -38996
$ws.Range("M80").Value = ""

$ws.Range("H83").Value = 37000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 37000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 111000
$ws.Range("N83").Value = -120984
$ws.Range("M83").Value = ""

$ws.Range("H97").Value = 577.5
$ws.Range("I97").Value = 448.8889
$ws.Range("J97").Value = 963.3333
$ws.Range("K97").Value = 448.8889
$ws.Range("L97").Value = 963.3333
$ws.Range("M97").Value = 47.11110000000002
$ws.Range("N97").Value = -1955.3333

$ws.Range("H136").Value = 950.2632
$ws.Range("I136").Value = 776.6875
$ws.Range("J136").Value = 1876
$ws.Range("K136").Value = 2330.0625
$ws.Range("L136").Value = 5628
$ws.Range("M136").Value = 219.9375
$ws.Range("N136").Value = -10728

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = ""

$ws.Range("H94").Value = 62501176
$ws.Range("I94").Value = 125000350
$ws.Range("J94").Value = 1999.5
$ws.Range("K94").Value = 125000350
$ws.Range("L94").Value = 1999.5
$ws.Range("M94").Value = -124999899
$ws.Range("N94").Value = -2901.5

$ws.Range("H98").Value = 59000
$ws.Range("I98").Value = 59000
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 59000
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -56005
$ws.Range("N98").Value = ""

$ws.Range("H140").Value = 22721.176
$ws.Range("I140").Value = 20780
$ws.Range("J140").Value = 22842.5
$ws.Range("K140").Value = 20780
$ws.Range("L140").Value = 22842.5
$ws.Range("M140").Value = -15600
$ws.Range("N140").Value = -33202.5

$ws.Range("H141").Value = 89993.336
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 89993.336
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 89993.336
$ws.Range("N141").Value = -100353.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4300
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 4800
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 4800
$ws.Range("M80").Value = -802
$ws.Range("N80").Value = -6796

$ws.Range("H83").Value = 4300
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 4800
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 24000
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -33984

$ws.Range("H132").Value = 2651.04
$ws.Range("I132").Value = 2258.125
$ws.Range("J132").Value = 3349.5557
$ws.Range("K132").Value = 6774.375
$ws.Range("L132").Value = 10048.6671
$ws.Range("M132").Value = -4244.375
$ws.Range("N132").Value = -15108.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2363.3333
$ws.Range("I82").Value = 2245
$ws.Range("J82").Value = 2600
$ws.Range("K82").Value = 2245
$ws.Range("L82").Value = 2600
$ws.Range("M82").Value = -1884
$ws.Range("N82").Value = -3322

$ws.Range("H85").Value = 2363.3333
$ws.Range("I85").Value = 2245
$ws.Range("J85").Value = 2600
$ws.Range("K85").Value = 2245
$ws.Range("L85").Value = 2600
$ws.Range("M85").Value = -997
$ws.Range("N85").Value = -5096

$ws.Range("H93").Value = 1466.6666
$ws.Range("I93").Value = 1300
$ws.Range("J93").Value = 1800
$ws.Range("K93").Value = 1300
$ws.Range("L93").Value = 1800
$ws.Range("M93").Value = -52
$ws.Range("N93").Value = -4296

$ws.Range("H132").Value = 65000.375
$ws.Range("I132").Value = 2571.1428
$ws.Range("J132").Value = 113556.445
$ws.Range("K132").Value = 7713.428400000001
$ws.Range("L132").Value = 340669.335
$ws.Range("M132").Value = -5183.428400000001
$ws.Range("N132").Value = -345729.335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 822.1667
$ws.Range("I100").Value = 760.75
$ws.Range("J100").Value = 945
$ws.Range("K100").Value = 1521.5
$ws.Range("L100").Value = 1890
$ws.Range("M100").Value = -980.5
$ws.Range("N100").Value = -2972

$ws.Range("H107").Value = 639.8333
$ws.Range("I107").Value = 530
$ws.Range("J107").Value = 749.6667
$ws.Range("K107").Value = 1590
$ws.Range("L107").Value = 2249.0001
$ws.Range("M107").Value = 330
$ws.Range("N107").Value = -6089.0001

$ws.Range("H126").Value = 71430100
$ws.Range("I126").Value = 125000750
$ws.Range("J126").Value = 2557.5
$ws.Range("K126").Value = 375002250
$ws.Range("L126").Value = 7672.5
$ws.Range("M126").Value = -374999780
$ws.Range("N126").Value = -12612.5

$ws.Range("H132").Value = 4094.2222
$ws.Range("I132").Value = 5085.125
$ws.Range("J132").Value = 2652.9092
$ws.Range("K132").Value = 15255.375
$ws.Range("L132").Value = 7958.7276
$ws.Range("M132").Value = -12725.375
$ws.Range("N132").Value = -13018.7276

$ws.Range("H136").Value = 459.9524
$ws.Range("I136").Value = 386.76923
$ws.Range("J136").Value = 578.875
$ws.Range("K136").Value = 1160.30769
$ws.Range("L136").Value = 1736.625
$ws.Range("M136").Value = 1389.69231
$ws.Range("N136").Value = -6836.625
